$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy the formatting from column K into column L for rows 4-13 ---
# (column L is brand-new content for most of these rows, so pull in the
# neighbouring cell's style first, then set the value)
$ws.Range("K4:K13").Copy()
$ws.Range("L4:L13").PasteSpecial(-4122)

# --- Populate the new "2020" column (L) ---
$ws.Range("L4").Value = 2020
$ws.Range("L5").Value = 1.2
$ws.Range("L6").Value = 1.7
$ws.Range("L7").Value = 0.4
$ws.Range("L8").Value = 3.3
$ws.Range("L9").Value = 3.9
$ws.Range("L10").Value = 2.4
$ws.Range("L11").Value = 95.5
$ws.Range("L12").Value = 94.4
$ws.Range("L13").Value = 97.2

# --- Update the view: scroll so column C is the left-most visible column,
#     and select the newly added L4:L13 block (active cell L4) ---
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("L4:L13").Select()
